$d = $word.ActiveDocument

# --- Occurrence #1: "I read the Online Secured Cash Bond Posting Summary (form CR-233)..." ---
# Replace just the title text (leave the following space/"(" run alone for now).
$rng1 = $d.Content
$rng1.Find.Execute("Online Secured Cash Bond Posting Summary", $false, $false, $false, $false, $false, `
                    $true, 1, $false, "Notice of Request to Post Secured Cash Bail Online", 1)

# The single space right after the new title (before "(") becomes italic, splitting it
# off from the "(form CR-233) ..." run, matching the author's fine-grained formatting tweak.
$spaceStart = $rng1.End
$spaceRng = $d.Range($spaceStart, $spaceStart + 1)
$spaceRng.Italic = 1

# --- Occurrence #2: "...defendant's CR-233 Online Secured Cash Bond Posting Summary for court case..." ---
$rng2 = $d.Content
$rng2.Find.Execute("Online Secured Cash Bond Posting Summary", $false, $false, $false, $false, $false, `
                    $true, 1, $false, "Notice of Request to Post Secured Cash Bail Online", 1)
